# Applies updated voltage-magnitude (vm_pu) results for the 380 kV case
# (commit: "case with 380 kV done") to Sheet1 of the active workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = [ordered]@{
    "B2" = 1.02
    "C2" = 1.042616957467731
    "D2" = 1.05096233040263
    "E2" = 1.051091100870098
    "F2" = 1.06279069561123
    "I2" = 1.042499088653712
    "J2" = 1.047691817825404
    "K2" = 1.053715152177703
    "L2" = 1.053843565336549
    "M2" = 1.065511087553805
    "N2" = 1.049179660378648
    "B3" = 1.02
    "C3" = 1.043611792977963
    "D3" = 1.051745002603919
    "E3" = 1.051956871593313
    "F3" = 1.063698177946945
    "I3" = 1.042720277948297
    "J3" = 1.048332939573941
    "K3" = 1.054310452525342
    "L3" = 1.054521776084281
    "M3" = 1.066233221930496
    "N3" = 1.049821692593605
    "B4" = 1.02
    "C4" = 1.044255879456775
    "D4" = 1.052251428599154
    "E4" = 1.052517692137305
    "F4" = 1.064285796852853
    "I4" = 1.042861779289435
    "J4" = 1.048747541229604
    "K4" = 1.054694978086035
    "L4" = 1.054960590639244
    "M4" = 1.066700273140494
    "N4" = 1.050236883031174
    "B5" = 1.02
    "C5" = 1.044526739203339
    "D5" = 1.052464325179201
    "E5" = 1.052753605690222
    "F5" = 1.064532930426142
    "I5" = 1.042920877305246
    "J5" = 1.048921779784497
    "K5" = 1.054856470399931
    "L5" = 1.055145059386934
    "M5" = 1.066896568183423
    "N5" = 1.0504113690248
    "B6" = 1.02
    "C6" = 1.044572222756435
    "D6" = 1.052500071127355
    "E6" = 1.05279322505815
    "F6" = 1.064574430984301
    "I6" = 1.042930777284341
    "J6" = 1.048951031662428
    "K6" = 1.054883576128206
    "L6" = 1.055176031952186
    "M6" = 1.066929523828388
    "N6" = 1.050440662443754
    "B7" = 1.02
    "C7" = 1.044259498361771
    "D7" = 1.052254273354551
    "E7" = 1.052520843859091
    "F7" = 1.064289098678
    "I7" = 1.04286257048988
    "J7" = 1.048749869650462
    "K7" = 1.054697136592933
    "L7" = 1.05496305555582
    "M7" = 1.06670289625522
    "N7" = 1.050239214758657
    "B8" = 1.02
    "C8" = 1.042953091524062
    "D8" = 1.051226840542107
    "E8" = 1.051383565323846
    "F8" = 1.063097296371304
    "I8" = 1.04257417627818
    "J8" = 1.047908538550443
    "K8" = 1.05391647559796
    "L8" = 1.054072776124363
    "M8" = 1.065755180734573
    "N8" = 1.049396688871974
    "B9" = 1.02
    "C9" = 1.040653835241551
    "D9" = 1.049416310980646
    "E9" = 1.049384252070489
    "F9" = 1.061000439652757
    "I9" = 1.042053586210278
    "J9" = 1.046424151428842
    "K9" = 1.052535735309763
    "L9" = 1.052503778077972
    "M9" = 1.064083563044362
    "N9" = 1.047910193750398
    "B10" = 1.02
    "C10" = 1.039122920699602
    "D10" = 1.048209322865195
    "E10" = 1.048054617659269
    "F10" = 1.059604795635821
    "I10" = 1.041698222977736
    "J10" = 1.045433362499494
    "K10" = 1.05161185872137
    "L10" = 1.051457692253817
    "M10" = 1.062968124097574
    "N10" = 1.046917997787126
    "B11" = 1.02
    "C11" = 1.038460481848157
    "D11" = 1.047686706172115
    "E11" = 1.047479654341878
    "F11" = 1.059001017790879
    "I11" = 1.041542384609242
    "J11" = 1.04500406588519
    "K11" = 1.051211019066496
    "L11" = 1.051004716175836
    "M11" = 1.062484896020963
    "N11" = 1.046488091522388
    "B12" = 1.02
    "C12" = 1.0382144918211
    "D12" = 1.047492586795591
    "E12" = 1.047266205237373
    "F12" = 1.058776830911547
    "I12" = 1.0414842048044
    "J12" = 1.044844564922972
    "K12" = 1.05106201080654
    "L12" = 1.050836459450887
    "M12" = 1.062305369183447
    "N12" = 1.046328364050509
    "B13" = 1.02
    "C13" = 1.038267254373859
    "D13" = 1.047534225873311
    "E13" = 1.04731198544307
    "F13" = 1.058824915978995
    "I13" = 1.041496697888043
    "J13" = 1.044878780267944
    "K13" = 1.051093978935416
    "L13" = 1.050872551131838
    "M13" = 1.062343879840958
    "N13" = 1.046362627985196
    "B14" = 1.02
    "C14" = 1.038440146819698
    "D14" = 1.047670660121054
    "E14" = 1.047462008161337
    "F14" = 1.058982484731551
    "I14" = 1.041537581452971
    "J14" = 1.044990882319679
    "K14" = 1.051198704410882
    "L14" = 1.050990808042783
    "M14" = 1.062470056972892
    "N14" = 1.046474889234703
    "B15" = 1.02
    "C15" = 1.038546680755864
    "D15" = 1.047754722293246
    "E15" = 1.047554457761843
    "F15" = 1.059079579086685
    "I15" = 1.041562732159016
    "J15" = 1.045059946671861
    "K15" = 1.051263213547823
    "L15" = 1.051063669889995
    "M15" = 1.06254779433761
    "N15" = 1.046544051666186
    "B16" = 1.02
    "C16" = 1.039166894266015
    "D16" = 1.04824400767383
    "E16" = 1.048092792549888
    "F16" = 1.059644877971077
    "I16" = 1.041708524128824
    "J16" = 1.045461847699086
    "K16" = 1.051638444454183
    "L16" = 1.051487754562863
    "M16" = 1.063000189473185
    "N16" = 1.046946523438969
    "B17" = 1.02
    "C17" = 1.039556060628911
    "D17" = 1.048550928871465
    "E17" = 1.048430684593142
    "F17" = 1.059999621833469
    "I17" = 1.041799450207451
    "J17" = 1.045713875572269
    "K17" = 1.051873604663575
    "L17" = 1.051753768533197
    "M17" = 1.063283902598619
    "N17" = 1.04719890922064
    "B18" = 1.02
    "C18" = 1.039783098864258
    "D18" = 1.048729952388188
    "E18" = 1.048627846215178
    "F18" = 1.060206590453766
    "I18" = 1.041852296342241
    "J18" = 1.045860852248103
    "K18" = 1.052010692871108
    "L18" = 1.051908928571571
    "M18" = 1.063449364861416
    "N18" = 1.047346094620213
    "B19" = 1.02
    "C19" = 1.039860520495014
    "D19" = 1.048790995012783
    "E19" = 1.048695085905305
    "F19" = 1.06027717027826
    "I19" = 1.041870283348326
    "J19" = 1.045910962907593
    "K19" = 1.052057423360052
    "L19" = 1.051961833869488
    "M19" = 1.063505779319125
    "N19" = 1.047396276442586
    "B20" = 1.02
    "C20" = 1.03951430219767
    "D20" = 1.048517998960308
    "E20" = 1.048394424218313
    "F20" = 1.059961555732575
    "I20" = 1.041789714297594
    "J20" = 1.045686838145484
    "K20" = 1.051848382132337
    "L20" = 1.051725227891255
    "M20" = 1.063253465216736
    "N20" = 1.047171833397608
    "B21" = 1.02
    "C21" = 1.038389232387267
    "D21" = 1.047630483523288
    "E21" = 1.04741782694726
    "F21" = 1.058936082357901
    "I21" = 1.041525550387346
    "J21" = 1.044957872187635
    "K21" = 1.051167868634529
    "L21" = 1.050955984357304
    "M21" = 1.06243290189448
    "N21" = 1.046441832224485
    "B22" = 1.02
    "C22" = 1.037682257325402
    "D22" = 1.047072489510899
    "E22" = 1.046804484001171
    "F22" = 1.05829180789423
    "I22" = 1.041357756452176
    "J22" = 1.044499305022175
    "K22" = 1.050739317240852
    "L22" = 1.050472323664052
    "M22" = 1.061916782544428
    "N22" = 1.045982613841052
    "B23" = 1.02
    "C23" = 1.038056999985714
    "D23" = 1.047368290206895
    "E23" = 1.047129563635816
    "F23" = 1.058633303957429
    "I23" = 1.041446868495636
    "J23" = 1.044742422398927
    "K23" = 1.050966565133895
    "L23" = 1.050728721700543
    "M23" = 1.062190405674327
    "N23" = 1.046226076472365
    "B24" = 1.02
    "C24" = 1.039533170899048
    "D24" = 1.048532878564926
    "E24" = 1.04841080848951
    "F24" = 1.059978756002814
    "I24" = 1.041794114121754
    "J24" = 1.045699055278041
    "K24" = 1.051859779345747
    "L24" = 1.051738124183237
    "M24" = 1.063267218629964
    "N24" = 1.047184067879894
    "B25" = 1.02
    "C25" = 1.04124791215469
    "D25" = 1.049884376034758
    "E25" = 1.049900556302729
    "F25" = 1.061542134397958
    "I25" = 1.042189637315633
    "J25" = 1.046808115743321
    "K25" = 1.052893290187217
    "L25" = 1.052909421083934
    "M25" = 1.064515901950715
    "N25" = 1.048294703338246
}

foreach ($cellRef in $newValues.Keys) {
    $ws.Range($cellRef).Value = $newValues[$cellRef]
}
